$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by writing values; build full data array for rows 2-31
$data = @(
    @('Entiat River Lake 03', 'Entiat River-Lake Creek', 'Entiat', 'HQ_spring_chinook', 1, 'NA', 0, 'Habitat Protection'),
    @('Entiat River Lake 05', 'Entiat River-Lake Creek', 'Entiat', 'HQ_spring_chinook', 1, 'NA', 0, 'Habitat Protection'),
    @('Entiat River Lake 06', 'Entiat River-Lake Creek', 'Entiat', 'HQ_spring_chinook,LF_spring_chinook', 2, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Entiat River Potato 05', 'Entiat River-Potato Creek', 'Entiat', 'HQ_spring_chinook,HQ_steelhead', 2, 'NA', 0, 'Habitat Protection'),
    @('Methow River Fawn 04', 'Methow River-Fawn Creek', 'Methow', 'HQ_spring_chinook,HQ_steelhead', 2, 'NA', 0, 'Habitat Protection'),
    @('Methow River Rattlesnake 06', 'Methow River-Rattlesnake Creek', 'Methow', 'HQ_spring_chinook', 1, 'NA', 0, 'Habitat Protection'),
    @('Nason Creek Lower 01', 'Lower Nason Creek', 'Wenatchee', 'HQ_spring_chinook,HQ_steelhead,LF_spring_chinook,LF_steelhead', 4, 'Spawning and Incubation,Winter Rearing,Summer Rearing', 3, 'Habitat Protection'),
    @('Nason Creek Lower 02', 'Lower Nason Creek', 'Wenatchee', 'HQ_spring_chinook,HQ_steelhead,LF_spring_chinook,LF_steelhead', 4, 'Summer Rearing,Spawning and Incubation', 2, 'Habitat Protection'),
    @('Nason Creek Lower 03', 'Lower Nason Creek', 'Wenatchee', 'HQ_spring_chinook,HQ_steelhead,LF_spring_chinook,LF_steelhead', 4, 'Spawning and Incubation,Winter Rearing,Summer Rearing', 3, 'Habitat Protection'),
    @('Twisp River Middle 06', 'Middle Twisp River', 'Methow', 'HQ_spring_chinook,HQ_steelhead', 2, 'NA', 0, 'Habitat Protection'),
    @('Little Wenatchee River Lower 05', 'Lower Little Wenatchee River', 'Wenatchee', 'LF_spring_chinook', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('White River Lower 08', 'Lower White River', 'Wenatchee', 'LF_spring_chinook,LF_bull_trout', 2, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Entiat River Potato 08', 'Entiat River-Potato Creek', 'Entiat', 'LF_spring_chinook,LF_steelhead', 2, 'Winter Rearing,Summer Rearing', 2, 'Habitat Protection'),
    @('Entiat River Preston 02', 'Entiat River-Preston Creek', 'Entiat', 'LF_spring_chinook', 1, 'Winter Rearing,Summer Rearing', 2, 'Habitat Protection'),
    @('Methow River Thompson 07', 'Methow River-Thompson Creek', 'Methow', 'LF_spring_chinook,LF_steelhead', 2, 'Winter Rearing,Summer Rearing', 2, 'Habitat Protection'),
    @('Entiat River Lake 01', 'Entiat River-Lake Creek', 'Entiat', 'LF_spring_chinook', 1, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Entiat River Lake 07', 'Entiat River-Lake Creek', 'Entiat', 'LF_spring_chinook', 1, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Entiat River Lake 08', 'Entiat River-Lake Creek', 'Entiat', 'LF_spring_chinook', 1, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Entiat River Lake 10', 'Entiat River-Lake Creek', 'Entiat', 'LF_spring_chinook', 1, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Entiat River Preston 05', 'Entiat River-Lake Creek', 'Entiat', 'LF_spring_chinook', 1, 'Summer Rearing', 1, 'Habitat Protection'),
    @('Nason Creek Lower 05', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 06', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 07', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 08', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 10', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 12', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 13', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 14', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Nason Creek Lower 15', 'Lower Nason Creek', 'Wenatchee', 'LF_steelhead', 1, 'Spawning and Incubation', 1, 'Habitat Protection'),
    @('Entiat River Mills 06', 'Entiat River-Mills Creek', 'Entiat', 'LF_steelhead', 1, 'Winter Rearing', 1, 'Habitat Protection')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$wb.Save()